# Auto-generated script to apply cryptos.xlsx price/volume updates
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, [string]$value) {
    # Force the cell to be treated as text so Excel does not
    # auto-convert numeric-looking strings (e.g. "1.00") into numbers,
    # and restore the default (unstyled) cell style afterwards.
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

Set-TextValue $ws.Range("D2") "38.095.42"
Set-TextValue $ws.Range("E2") "  +1.78%  "

Set-TextValue $ws.Range("D3") "2.052.52"
Set-TextValue $ws.Range("E3") "  +0.84%  "

Set-TextValue $ws.Range("E4") "  -0.04%  "

Set-TextValue $ws.Range("D5") "228.16"
Set-TextValue $ws.Range("E5") "  -0.83%  "

Set-TextValue $ws.Range("E6") "  +0.08%  "

Set-TextValue $ws.Range("D7") "61.10"
Set-TextValue $ws.Range("E7") "  +8.36%  "

Set-TextValue $ws.Range("E8") "  -0.01%  "

Set-TextValue $ws.Range("E9") "  +0.33%  "

Set-TextValue $ws.Range("D10") "0.0819"
Set-TextValue $ws.Range("E10") "  +2.68%  "

Set-TextValue $ws.Range("E11") "  +0.92%  "

Set-TextValue $ws.Range("D12") "14.77"
Set-TextValue $ws.Range("E12") "  +2.19%  "

Set-TextValue $ws.Range("D13") "2.356.43"
Set-TextValue $ws.Range("E13") "  +0.91%  "

Set-TextValue $ws.Range("D14") "21.08"
Set-TextValue $ws.Range("E14") "  +3.36%  "

Set-TextValue $ws.Range("D15") "0.760"
Set-TextValue $ws.Range("E15") "  +2.04%  "

Set-TextValue $ws.Range("D16") "5.29"
Set-TextValue $ws.Range("E16") "  +0.96%  "

Set-TextValue $ws.Range("D17") "2.041.68"
Set-TextValue $ws.Range("E17") "  +0.26%  "

Set-TextValue $ws.Range("D18") "38.036.94"
Set-TextValue $ws.Range("E18") "  +1.80%  "

Set-TextValue $ws.Range("E19") "  -1.53%  "

Set-TextValue $ws.Range("D20") "69.81"
Set-TextValue $ws.Range("E20") "  +1.04%  "

Set-TextValue $ws.Range("D21") "0.0₃0829"
Set-TextValue $ws.Range("E21") "  +0.46%  "

Set-TextValue $ws.Range("D22") "224.78"
Set-TextValue $ws.Range("E22") "  +0.39%  "

Set-TextValue $ws.Range("D23") "1.00"
Set-TextValue $ws.Range("E23") "  +0.00%  "

Set-TextValue $ws.Range("E24") "  -0.85%  "

Set-TextValue $ws.Range("D25") "2.21"
Set-TextValue $ws.Range("E25") "  -1.45%  "

Set-TextValue $ws.Range("D26") "166.48"
Set-TextValue $ws.Range("E26") "  +0.98%  "

Set-TextValue $ws.Range("D27") "9.22"
Set-TextValue $ws.Range("E27") "  +0.35%  "

Set-TextValue $ws.Range("D28") "0.131"
Set-TextValue $ws.Range("E28") "  -1.55%  "

Set-TextValue $ws.Range("D29") "18.97"
Set-TextValue $ws.Range("E29") "  +0.87%  "

Set-TextValue $ws.Range("E30") "  -2.58%  "

Set-TextValue $ws.Range("E31") "  +2.25%  "

Set-TextValue $ws.Range("D32") "4.49"
Set-TextValue $ws.Range("E32") "  -0.17%  "

Set-TextValue $ws.Range("E33") "  +1.35%  "

Set-TextValue $ws.Range("D34") "4.53"
Set-TextValue $ws.Range("E34") "  +1.08%  "

Set-TextValue $ws.Range("E35") "  -0.57%  "

Set-TextValue $ws.Range("D36") "6.31"
Set-TextValue $ws.Range("E36") "  +9.58%  "

Set-TextValue $ws.Range("E37") "  -1.85%  "

Set-TextValue $ws.Range("E38") "  +0.65%  "

Set-TextValue $ws.Range("E39") "  +0.15%  "

Set-TextValue $ws.Range("D40") "1.535.72"
Set-TextValue $ws.Range("E40") "  +4.47%  "

Set-TextValue $ws.Range("B41") "Aave"
Set-TextValue $ws.Range("C41") "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
Set-TextValue $ws.Range("D41") "97.73"
Set-TextValue $ws.Range("E41") "  +3.11%  "

Set-TextValue $ws.Range("B42") "VeChain"
Set-TextValue $ws.Range("C42") "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
Set-TextValue $ws.Range("D42") "0.0217"
Set-TextValue $ws.Range("E42") "  +1.45%  "

Set-TextValue $ws.Range("B43") "InjectiveProtocol"
Set-TextValue $ws.Range("C43") "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
Set-TextValue $ws.Range("D43") "16.66"
Set-TextValue $ws.Range("E43") "  +1.97%  "

Set-TextValue $ws.Range("B44") "HuobiToken"
Set-TextValue $ws.Range("C44") "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
Set-TextValue $ws.Range("D44") "2.84"
Set-TextValue $ws.Range("E44") "  +0.95%  "

Set-TextValue $ws.Range("D45") "0.0928"
Set-TextValue $ws.Range("E45") "  -0.29%  "

Set-TextValue $ws.Range("E46") "  +0.71%  "

Set-TextValue $ws.Range("D47") "4.00"
Set-TextValue $ws.Range("E47") "  -5.91%  "

Set-TextValue $ws.Range("E48") "  +1.93%  "

Set-TextValue $ws.Range("E49") "  -0.08%  "

Set-TextValue $ws.Range("D50") "7.07"
Set-TextValue $ws.Range("E50") "  -0.78%  "

Set-TextValue $ws.Range("D51") "2.246.01"
Set-TextValue $ws.Range("E51") "  +1.04%  "

